$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = "[49.96419634525922, 50.05653482258821]"
$ws.Range("T2").Value = "[49.96085291858227, 50.02668313223217]"
$ws.Range("L3").Value = "[49.96445084161835, 50.05969254209422]"
$ws.Range("T3").Value = "[49.96908925346484, 50.02587327650247]"
